$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text representation
# (values such as "322.05" or "48.827.02" must stay as literal strings,
# not be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '48.827.02'
$ws.Range('E2').Value = '  -2.16%  '
$ws.Range('D3').Value = '2.614.88'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').Value = '322.05'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('E7').Value = '  -1.53%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -4.04%  '
$ws.Range('D10').Value = '39.40'
$ws.Range('E10').Value = '  -3.47%  '
$ws.Range('D11').Value = '19.64'
$ws.Range('E11').Value = '  -5.42%  '
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').Value = '7.20'
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').Value = '3.025.55'
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').Value = '2.608.44'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '48.759.23'
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('D19').Value = '2.97'
$ws.Range('E19').Value = '  -3.50%  '
$ws.Range('D20').Value = '12.83'
$ws.Range('E20').Value = '  -4.14%  '
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('D22').Value = '0.0₃0943'
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('D23').Value = '268.43'
$ws.Range('E23').Value = '  -5.25%  '
$ws.Range('D24').Value = '68.64'
$ws.Range('E24').Value = '  -5.78%  '
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('D26').Value = '25.99'
$ws.Range('E26').Value = '  -2.56%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = '10.01'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('E30').Value = '  -3.01%  '
$ws.Range('E31').Value = '  -7.16%  '
$ws.Range('D32').Value = '49.37'
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('D33').Value = '5.48'
$ws.Range('E33').Value = '  +0.69%  '
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').Value = '0.0795'
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '5.01'
$ws.Range('E36').Value = '  +5.56%  '
$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').Value = '18.86'
$ws.Range('E37').Value = '  -4.47%  '
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').Value = '3.10'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('D40').Value = '126.16'
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('E41').Value = '  -1.59%  '
$ws.Range('D42').Value = '22.15'
$ws.Range('E42').Value = '  -4.52%  '
$ws.Range('D43').Value = '2.12'
$ws.Range('E43').Value = '  -4.36%  '
$ws.Range('D44').Value = '0.0317'
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').Value = '2.058.09'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('E46').Value = '  -3.67%  '
$ws.Range('D47').Value = '2.13'
$ws.Range('E47').Value = '  +5.61%  '
$ws.Range('D48').Value = '2.15'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('D49').Value = '8.85'
$ws.Range('E49').Value = '  -3.44%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '5.15'
$ws.Range('E50').Value = '  -4.22%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = '58.31'
$ws.Range('E51').Value = '  +1.33%  '
